$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert two new paragraphs right after the existing blank paragraph
#    that follows the "DATE: ..." paragraph, and before the
#    "PLAINTIFF/TENANT: ..." paragraph:
#      - a new blank paragraph (same tab stops / rPr as its neighbours)
#      - a new "COURT LOCATION: {{ trial_court.address.on_one_line() }}"
#        paragraph
# ---------------------------------------------------------------------

$datePara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "^DATE:\s+The ") {
        $datePara = $cand
        break
    }
}

$blankPara = $datePara.Next()

# Create the two new paragraphs after the blank paragraph that already
# follows the DATE paragraph.
$blankPara.Range.InsertParagraphAfter()
$newBlank = $blankPara.Next()
$newBlank.Range.InsertParagraphAfter()
$newCourt = $newBlank.Next()

$tabsPPr = @'
<w:pPr>
  <w:tabs>
    <w:tab w:val="left" w:pos="540"/>
    <w:tab w:val="left" w:pos="720"/>
    <w:tab w:val="left" w:pos="1440"/>
    <w:tab w:val="left" w:pos="2160"/>
    <w:tab w:val="left" w:pos="2880"/>
    <w:tab w:val="left" w:pos="3600"/>
    <w:tab w:val="left" w:pos="4320"/>
    <w:tab w:val="left" w:pos="4500"/>
    <w:tab w:val="left" w:pos="4860"/>
    <w:tab w:val="left" w:pos="5040"/>
    <w:tab w:val="left" w:pos="5580"/>
  </w:tabs>
  <w:rPr>
    <w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/>
  </w:rPr>
</w:pPr>
'@

# The plain blank paragraph: pPr only, no runs.
$blankXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
$tabsPPr
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$newBlank.Range.InsertXML($blankXml)

# The "COURT LOCATION: {{ trial_court.address.on_one_line() }}" paragraph,
# split into runs the way Word's proofing tools would (matches the
# author's original edit exactly).
$courtXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
$tabsPPr
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t xml:space="preserve">COURT LOCATION: </w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t xml:space="preserve">{{ </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>trial</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>_</w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>court.address</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>.on_one_</w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>line</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>()</w:t>
</w:r>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t xml:space="preserve"> }</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
  <w:rPr><w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/></w:rPr>
  <w:t>}</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$newCourt.Range.InsertXML($courtXml)

# ---------------------------------------------------------------------
# 2) Un-hide the "Default Paragraph Font" character style (remove the
#    <w:semiHidden/> flag it carried).
# ---------------------------------------------------------------------
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $false

Write-Output "edit complete"
